$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 39.66867433333334
$ws.Range("H2").Value = 119.006023
$ws.Range("I2").Value = 0.154574216411057
$ws.Range("J2").Value = 0.1545742164110569
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 2.041175666666667
$ws.Range("N2").Value = 6.123527
$ws.Range("O2").Value = 0.1007574239555886
$ws.Range("P2").Value = 0.1007574239555885
$ws.Range("Q2").Value = 80.97073277812456
$ws.Range("R2").Value = 728.736595003121
$ws.Range("S2").Value = 0.01557449985553176
$ws.Range("T2").Value = 0.01557449985553175

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 39.66867433333334
$ws.Range("H3").Value = 119.006023
$ws.Range("I3").Value = 0.154574216411057
$ws.Range("J3").Value = 0.1545742164110569
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 9.510473333333332
$ws.Range("N3").Value = 28.53142
$ws.Range("O3").Value = 0.4694602279037812
$ws.Range("P3").Value = 0.4694602279037812
$ws.Range("Q3").Value = 377.2678694158511
$ws.Range("R3").Value = 3395.41082474266
$ws.Range("S3").Value = 0.0725664468643832
$ws.Range("T3").Value = 0.07256644686438317

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 39.66867433333334
$ws.Range("H4").Value = 119.006023
$ws.Range("I4").Value = 0.154574216411057
$ws.Range("J4").Value = 0.1545742164110569
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 6.660130333333334
$ws.Range("N4").Value = 19.980391
$ws.Range("O4").Value = 0.3287603250194579
$ws.Range("P4").Value = 0.3287603250194578
$ws.Range("Q4").Value = 264.1985412105548
$ws.Range("R4").Value = 2377.786870894993
$ws.Range("S4").Value = 0.05081786962692711
$ws.Range("T4").Value = 0.05081786962692709

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 39.66867433333334
$ws.Range("H5").Value = 119.006023
$ws.Range("I5").Value = 0.154574216411057
$ws.Range("J5").Value = 0.1545742164110569
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 2.046536
$ws.Range("N5").Value = 6.139608
$ws.Range("O5").Value = 0.1010220231211724
$ws.Range("P5").Value = 0.1010220231211723
$ws.Range("Q5").Value = 81.18337009544268
$ws.Range("R5").Value = 730.650330858984
$ws.Range("S5").Value = 0.01561540006421489
$ws.Range("T5").Value = 0.01561540006421489

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 57.66057933333332
$ws.Range("H6").Value = 172.981738
$ws.Range("I6").Value = 0.2246820449144221
$ws.Range("J6").Value = 0.2246820449144221
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 2.041175666666667
$ws.Range("N6").Value = 6.123527
$ws.Range("O6").Value = 0.1007574239555886
$ws.Range("P6").Value = 0.1007574239555885
$ws.Range("Q6").Value = 117.6953714611029
$ws.Range("R6").Value = 1059.258343149926
$ws.Range("S6").Value = 0.02263838405465101
$ws.Range("T6").Value = 0.02263838405465101

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 57.66057933333332
$ws.Range("H7").Value = 172.981738
$ws.Range("I7").Value = 0.2246820449144221
$ws.Range("J7").Value = 0.2246820449144221
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 9.510473333333332
$ws.Range("N7").Value = 28.53142
$ws.Range("O7").Value = 0.4694602279037812
$ws.Range("P7").Value = 0.4694602279037812
$ws.Range("Q7").Value = 548.3794021342176
$ws.Range("R7").Value = 4935.414619207959
$ws.Range("S7").Value = 0.1054792840114122
$ws.Range("T7").Value = 0.1054792840114122

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 57.66057933333332
$ws.Range("H8").Value = 172.981738
$ws.Range("I8").Value = 0.2246820449144221
$ws.Range("J8").Value = 0.2246820449144221
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 6.660130333333334
$ws.Range("N8").Value = 19.980391
$ws.Range("O8").Value = 0.3287603250194579
$ws.Range("P8").Value = 0.3287603250194578
$ws.Range("Q8").Value = 384.0269734555064
$ws.Range("R8").Value = 3456.242761099558
$ws.Range("S8").Value = 0.07386654211210185
$ws.Range("T8").Value = 0.07386654211210182

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 57.66057933333332
$ws.Range("H9").Value = 172.981738
$ws.Range("I9").Value = 0.2246820449144221
$ws.Range("J9").Value = 0.2246820449144221
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 2.046536
$ws.Range("N9").Value = 6.139608
$ws.Range("O9").Value = 0.1010220231211724
$ws.Range("P9").Value = 0.1010220231211723
$ws.Range("Q9").Value = 118.0044513865227
$ws.Range("R9").Value = 1062.040062478704
$ws.Range("S9").Value = 0.02269783473625703
$ws.Range("T9").Value = 0.02269783473625702

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 64.993678
$ws.Range("H10").Value = 194.981034
$ws.Range("I10").Value = 0.2532564301015895
$ws.Range("J10").Value = 0.2532564301015895
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 2.041175666666667
$ws.Range("N10").Value = 6.123527
$ws.Range("O10").Value = 0.1007574239555886
$ws.Range("P10").Value = 0.1007574239555885
$ws.Range("Q10").Value = 132.6635140207687
$ws.Range("R10").Value = 1193.971626186918
$ws.Range("S10").Value = 0.02551746549722473
$ws.Range("T10").Value = 0.02551746549722473

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 64.993678
$ws.Range("H11").Value = 194.981034
$ws.Range("I11").Value = 0.2532564301015895
$ws.Range("J11").Value = 0.2532564301015895
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 9.510473333333332
$ws.Range("N11").Value = 28.53142
$ws.Range("O11").Value = 0.4694602279037812
$ws.Range("P11").Value = 0.4694602279037812
$ws.Range("Q11").Value = 618.1206414542532
$ws.Range("R11").Value = 5563.08577308828
$ws.Range("S11").Value = 0.1188938213935902
$ws.Range("T11").Value = 0.1188938213935902

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 64.993678
$ws.Range("H12").Value = 194.981034
$ws.Range("I12").Value = 0.2532564301015895
$ws.Range("J12").Value = 0.2532564301015895
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 6.660130333333334
$ws.Range("N12").Value = 19.980391
$ws.Range("O12").Value = 0.3287603250194579
$ws.Range("P12").Value = 0.3287603250194578
$ws.Range("Q12").Value = 432.8663663226994
$ws.Range("R12").Value = 3895.797296904294
$ws.Range("S12").Value = 0.08326066627346618
$ws.Range("T12").Value = 0.08326066627346616

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 64.993678
$ws.Range("H13").Value = 194.981034
$ws.Range("I13").Value = 0.2532564301015895
$ws.Range("J13").Value = 0.2532564301015895
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 2.046536
$ws.Range("N13").Value = 6.139608
$ws.Range("O13").Value = 0.1010220231211724
$ws.Range("P13").Value = 0.1010220231211723
$ws.Range("Q13").Value = 133.011901799408
$ws.Range("R13").Value = 1197.107116194672
$ws.Range("S13").Value = 0.02558447693730834
$ws.Range("T13").Value = 0.02558447693730834

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 94.308965
$ws.Range("H14").Value = 282.926895
$ws.Range("I14").Value = 0.3674873085729315
$ws.Range("J14").Value = 0.3674873085729314
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 2.041175666666667
$ws.Range("N14").Value = 6.123527
$ws.Range("O14").Value = 0.1007574239555886
$ws.Range("P14").Value = 0.1007574239555885
$ws.Range("Q14").Value = 192.5011645065183
$ws.Range("R14").Value = 1732.510480558665
$ws.Range("S14").Value = 0.03702707454818105
$ws.Range("T14").Value = 0.03702707454818104

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 94.308965
$ws.Range("H15").Value = 282.926895
$ws.Range("I15").Value = 0.3674873085729315
$ws.Range("J15").Value = 0.3674873085729314
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 9.510473333333332
$ws.Range("N15").Value = 28.53142
$ws.Range("O15").Value = 0.4694602279037812
$ws.Range("P15").Value = 0.4694602279037812
$ws.Range("Q15").Value = 896.9228967267666
$ws.Range("R15").Value = 8072.306070540899
$ws.Range("S15").Value = 0.1725206756343956
$ws.Range("T15").Value = 0.1725206756343955

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 94.308965
$ws.Range("H16").Value = 282.926895
$ws.Range("I16").Value = 0.3674873085729315
$ws.Range("J16").Value = 0.3674873085729314
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 6.660130333333334
$ws.Range("N16").Value = 19.980391
$ws.Range("O16").Value = 0.3287603250194579
$ws.Range("P16").Value = 0.3287603250194578
$ws.Range("Q16").Value = 628.1099985017718
$ws.Range("R16").Value = 5652.989986515945
$ws.Range("S16").Value = 0.1208152470069628
$ws.Range("T16").Value = 0.1208152470069627

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 94.308965
$ws.Range("H17").Value = 282.926895
$ws.Range("I17").Value = 0.3674873085729315
$ws.Range("J17").Value = 0.3674873085729314
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 2.046536
$ws.Range("N17").Value = 6.139608
$ws.Range("O17").Value = 0.1010220231211724
$ws.Range("P17").Value = 0.1010220231211723
$ws.Range("Q17").Value = 193.00669199524
$ws.Range("R17").Value = 1737.06022795716
$ws.Range("S17").Value = 0.03712431138339208
$ws.Range("T17").Value = 0.03712431138339207
